$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = "AnyTableProperty table of"
$ws.Range("C6").Select()
